$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2025-07-27 Sunday" "2025-07-28 Monday"

Replace-Text "775×8=" "636×4="
Replace-Text "960×3=" "438×4="
Replace-Text "701×8=" "885×5="
Replace-Text "705×8=" "595×7="
Replace-Text "201×5=" "720×5="
Replace-Text "607×3=" "563×3="
Replace-Text "121×4=" "969×8="
Replace-Text "573×3=" "267×3="
Replace-Text "266×7=" "413×7="
Replace-Text "655×6=" "269×4="
Replace-Text "636×5=" "290×8="
Replace-Text "366×6=" "327×5="
Replace-Text "516×4=" "247×7="
Replace-Text "799×6=" "125×2="
Replace-Text "780×4=" "897×9="
Replace-Text "346×3=" "304×9="
Replace-Text "160×5=" "509×4="
Replace-Text "658×9=" "180×7="
Replace-Text "426×3=" "886×8="
Replace-Text "416×5=" "346×6="
Replace-Text "933×3=" "836×8="
Replace-Text "640×5=" "887×4="
Replace-Text "184×3=" "174×7="
Replace-Text "202×9=" "993×8="
Replace-Text "929×7=" "741×6="
